$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: the test-data row now reuses the same email address as row 2
# (instead of the old "lgstester@gmail.com"), and the password value
# changes from the numeric-looking "123233" to "Test123". The leading
# apostrophes force these text entries to keep their existing
# text/quote-prefixed cell style instead of Excel re-typing them.
$ws.Range("A3").Value = "'lgstester50@gmail.com"
$ws.Range("B3").Value = "'Test123"

# Update the saved selection/active cell to B8.
$ws.Range("B8").Select()
